# Add two new event rows (137, 138) to the events table, matching the
# "RESONANZ" and "TECHNO NIGHT" entries at PM93 / Essen, each with an
# Instagram reel hyperlink in column E (styled like the other Link cells:
# single underline, green font).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$greenOle = 65280   # OLE (0x00BBGGRR) for RGB 00FF00 - matches the Link column's green underline styling
$noUnderline = -4142  # xlUnderlineStyleNone
$textFormat = "@"

function Set-LinkCellFormatting($range, [string]$url) {
    # Colour/underline the whole visible run of text (all-but-last char, then
    # the last char) so the shared string keeps one rich <r> run spanning the
    # full URL - matching how the workbook's other hyperlink cells are styled.
    $len = $url.Length
    $head = $range.Characters(1, $len - 1)
    $head.Font.Underline = 2
    $head.Font.Color = $greenOle
    $tail = $range.Characters($len, 1)
    $tail.Font.Underline = 2
    $tail.Font.Color = $greenOle
}

function Restore-PlainCellChrome($range) {
    # Hyperlinks.Add auto-applies Excel's built-in "Hyperlink" cell style;
    # put the direct formatting back to match the sheet's normal bordered /
    # text-formatted table cells (same look as the other populated rows).
    $range.NumberFormat = $textFormat
    $range.Interior.ColorIndex = 2
    $range.Font.Color = 0
    $range.Font.Underline = $noUnderline
}

function Add-EventRow([int]$row, [double]$date, [string]$event, [string]$location, [string]$city, [string]$url) {
    $ws.Range("A$row").Value = $date

    $ws.Range("B$row").NumberFormat = $textFormat
    $ws.Range("B$row").Value = $event

    $ws.Range("C$row").NumberFormat = $textFormat
    $ws.Range("C$row").Value = $location

    $ws.Range("D$row").NumberFormat = $textFormat
    $ws.Range("D$row").Value = $city

    $eRange = $ws.Range("E$row")
    $eRange.NumberFormat = $textFormat
    $eRange.Value = $url

    $ws.Hyperlinks.Add($eRange, $url, "", "", $url) | Out-Null

    Set-LinkCellFormatting $eRange $url
    Restore-PlainCellChrome $eRange
}

Add-EventRow 137 45710 "RESONANZ" "PM93" "Essen" "https://www.instagram.com/reel/DFqaRccOVWI/?igsh=MWI1anVrd21zcWVlaw=="
Add-EventRow 138 45696 "TECHNO NIGHT" "PM93" "Essen" "https://www.instagram.com/reel/DFqdAWZuoIC/?igsh=ajZrenlwNmFxNXc1"
